$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet ("comps")
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "spaceInColHeader"

# Populate header row - "col2 " intentionally has a trailing space to
# reproduce the whitespace-in-column-header bug being tested
$newSheet.Range("A1").Value = "col1"
$newSheet.Range("B1").Value = "col2 "

# Populate data rows
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 2
$newSheet.Range("A3").Value = 3
$newSheet.Range("B3").Value = 4

# Make the new sheet the active sheet/tab and select B1
$newSheet.Activate()
$newSheet.Range("B1").Select()
